# Auto-generated Excel COM-interop script
# Applies numeric value updates to the Leve profit-tracking sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) as produced by the scheduled
# market-price refresh run.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1553.7778
$ws.Range("I19").Value = 1839.8
$ws.Range("J19").Value = 1196.25
$ws.Range("K19").Value = 1839.8
$ws.Range("L19").Value = 1196.25
$ws.Range("M19").Value = -1664.8
$ws.Range("N19").Value = -1546.25
$ws.Range("H34").Value = 7957.2856
$ws.Range("I34").Value = 7957.2856
$ws.Range("K34").Value = 7957.2856
$ws.Range("M34").Value = -7754.2856
$ws.Range("H36").Value = 7957.2856
$ws.Range("I36").Value = 7957.2856
$ws.Range("K36").Value = 7957.2856
$ws.Range("M36").Value = -7242.2856
$ws.Range("H48").Value = 4000
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 4000
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 12000
$ws.Range("M48").ClearContents()
$ws.Range("N48").Value = -12584
$ws.Range("H56").Value = 4000
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 4000
$ws.Range("K56").Value = 0
$ws.Range("L56").Value = 12000
$ws.Range("M56").ClearContents()
$ws.Range("N56").Value = -13068
$ws.Range("H69").Value = 37177.92
$ws.Range("I69").Value = 117006.5
$ws.Range("K69").Value = 351019.5
$ws.Range("M69").Value = -350145.5
$ws.Range("H72").Value = 37177.92
$ws.Range("I72").Value = 117006.5
$ws.Range("K72").Value = 1053058.5
$ws.Range("M72").Value = -1048690.5
$ws.Range("H74").Value = 132065.67
$ws.Range("I74").Value = 141889.81
$ws.Range("K74").Value = 141889.81
$ws.Range("M74").Value = -140953.81
$ws.Range("H77").Value = 132065.67
$ws.Range("I77").Value = 141889.81
$ws.Range("K77").Value = 709449.05
$ws.Range("M77").Value = -704769.05
$ws.Range("H88").Value = 1831.5
$ws.Range("J88").Value = 1716.6666
$ws.Range("L88").Value = 1716.6666
$ws.Range("N88").Value = -2528.6666
$ws.Range("H91").Value = 1831.5
$ws.Range("J91").Value = 1716.6666
$ws.Range("L91").Value = 1716.6666
$ws.Range("N91").Value = -4524.6666
$ws.Range("H101").Value = 20003772
$ws.Range("I101").Value = 50008450
$ws.Range("J101").Value = 653
$ws.Range("K101").Value = 150025350
$ws.Range("L101").Value = 1959
$ws.Range("M101").Value = -150023728
$ws.Range("N101").Value = -5203
$ws.Range("H137").Value = 4276.4614
$ws.Range("I137").Value = 2288.5557
$ws.Range("K137").Value = 6865.6671
$ws.Range("M137").Value = -4315.6671
$ws.Range("H138").Value = 2472.1428
$ws.Range("J138").Value = 6000
$ws.Range("L138").Value = 18000
$ws.Range("N138").Value = -28280

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 30.9
$ws.Range("I5").Value = 30.9
$ws.Range("K5").Value = 30.9
$ws.Range("M5").Value = 81.09999999999999
$ws.Range("H29").Value = 13603
$ws.Range("I29").Value = 9009
$ws.Range("J29").Value = 15900
$ws.Range("K29").Value = 9009
$ws.Range("L29").Value = 15900
$ws.Range("M29").Value = -8701
$ws.Range("N29").Value = -16516
$ws.Range("H122").Value = 3721.1667
$ws.Range("I122").Value = 3332
$ws.Range("K122").Value = 9996
$ws.Range("M122").Value = -7546

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 30.9
$ws.Range("I4").Value = 30.9
$ws.Range("K4").Value = 30.9
$ws.Range("M4").Value = 84.09999999999999
$ws.Range("H86").Value = 1791.4286
$ws.Range("J86").Value = 1850
$ws.Range("L86").Value = 1850
$ws.Range("N86").Value = -4096
$ws.Range("H89").Value = 1791.4286
$ws.Range("J89").Value = 1850
$ws.Range("L89").Value = 9250
$ws.Range("N89").Value = -20482
$ws.Range("H94").Value = 1233
$ws.Range("I94").Value = 1233
$ws.Range("K94").Value = 1233
$ws.Range("M94").Value = -782
$ws.Range("H107").Value = 1933.4286
$ws.Range("I107").Value = 1933.4286
$ws.Range("K107").Value = 1933.4286
$ws.Range("M107").Value = -13.42859999999996

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("M42").ClearContents()
$ws.Range("H58").Value = 2349.2
$ws.Range("I58").Value = 1693.6666
$ws.Range("J58").Value = 3332.5
$ws.Range("K58").Value = 1693.6666
$ws.Range("L58").Value = 3332.5
$ws.Range("M58").Value = -1490.6666
$ws.Range("N58").Value = -3738.5
$ws.Range("H99").Value = 1894
$ws.Range("J99").Value = 2552
$ws.Range("L99").Value = 2552
$ws.Range("N99").Value = -5548
$ws.Range("H126").Value = 1894
$ws.Range("J126").Value = 2552
$ws.Range("L126").Value = 7656
$ws.Range("N126").Value = -12596
$ws.Range("H132").Value = 4199.2
$ws.Range("I132").Value = 3999.6667
$ws.Range("J132").Value = 4498.5
$ws.Range("K132").Value = 11999.0001
$ws.Range("L132").Value = 13495.5
$ws.Range("M132").Value = -9469.000100000001
$ws.Range("N132").Value = -18555.5
$ws.Range("H136").Value = 2349.2
$ws.Range("I136").Value = 1693.6666
$ws.Range("J136").Value = 3332.5
$ws.Range("K136").Value = 5080.9998
$ws.Range("L136").Value = 9997.5
$ws.Range("M136").Value = -2530.9998
$ws.Range("N136").Value = -15097.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 156
$ws.Range("I14").Value = 156
$ws.Range("K14").Value = 468
$ws.Range("M14").Value = -295
$ws.Range("H106").Value = 16466.334
$ws.Range("I106").Value = 9400
$ws.Range("J106").Value = 19999.5
$ws.Range("K106").Value = 28200
$ws.Range("L106").Value = 59998.5
$ws.Range("M106").Value = -27254
$ws.Range("N106").Value = -61890.5
$ws.Range("H121").Value = 14161.272
$ws.Range("I121").Value = 22055.2
$ws.Range("K121").Value = 66165.60000000001
$ws.Range("M121").Value = -64855.60000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4174.75
$ws.Range("I70").Value = 3999.6667
$ws.Range("K70").Value = 3999.6667
$ws.Range("M70").Value = -3729.6667
$ws.Range("H73").Value = 4174.75
$ws.Range("I73").Value = 3999.6667
$ws.Range("K73").Value = 3999.6667
$ws.Range("M73").Value = -3063.6667
$ws.Range("H126").Value = 7333.3335
$ws.Range("I126").Value = 5000
$ws.Range("K126").Value = 15000
$ws.Range("M126").Value = -12530

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1522.1305
$ws.Range("I22").Value = 1626.7894
$ws.Range("J22").Value = 1025
$ws.Range("K22").Value = 1626.7894
$ws.Range("L22").Value = 1025
$ws.Range("M22").Value = -1331.7894
$ws.Range("N22").Value = -1615
$ws.Range("H27").Value = 1522.1305
$ws.Range("I27").Value = 1626.7894
$ws.Range("J27").Value = 1025
$ws.Range("K27").Value = 1626.7894
$ws.Range("L27").Value = 1025
$ws.Range("M27").Value = -1519.7894
$ws.Range("N27").Value = -1239
$ws.Range("H132").Value = 3432.2222
$ws.Range("I132").Value = 3079.2
$ws.Range("J132").Value = 3873.5
$ws.Range("K132").Value = 9237.599999999999
$ws.Range("L132").Value = 11620.5
$ws.Range("M132").Value = -6707.599999999999
$ws.Range("N132").Value = -16680.5
$ws.Range("H136").Value = 3002
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 18444.5
$ws.Range("J41").Value = 17933.334
$ws.Range("L41").Value = 17933.334
$ws.Range("N41").Value = -18713.334
$ws.Range("H116").Value = 62341
$ws.Range("J116").Value = 62341
$ws.Range("L116").Value = 62341
$ws.Range("N116").Value = -71519
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()
$ws.Range("H132").Value = 3065.8
$ws.Range("I132").Value = 3065.8
$ws.Range("K132").Value = 9197.400000000001
$ws.Range("M132").Value = -6667.400000000001
